$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
# row 116 (diff @@ -6463)
$ws1.Range("H116").Value = 2139412.5
$ws1.Range("I116").Value = 2383656.8
$ws1.Range("J116").Value = 2275
$ws1.Range("K116").Value = 2383656.8
$ws1.Range("L116").Value = 2275
$ws1.Range("M116").Value = -2380214.8
$ws1.Range("N116").Value = -9159
# row 137 (diff @@ -7513)
$ws1.Range("H137").Value = 1574.3704
$ws1.Range("I137").Value = 1487.6875
$ws1.Range("J137").Value = 1700.4546
$ws1.Range("K137").Value = 4463.0625
$ws1.Range("L137").Value = 5101.3638
$ws1.Range("M137").Value = -1913.0625
$ws1.Range("N137").Value = -10201.3638

$ws2 = $wb.Worksheets.Item("ARM")
# row 2 (diff @@ -7861)
$ws2.Range("H2").Value = 630.30304
$ws2.Range("I2").Value = 550.3333
$ws2.Range("J2").Value = 696.94446
$ws2.Range("K2").Value = 550.3333
$ws2.Range("L2").Value = 696.94446
$ws2.Range("M2").Value = -437.3333
$ws2.Range("N2").Value = -922.94446
# row 32 (diff @@ -9319)
$ws2.Range("H32").Value = 15770.344
$ws2.Range("I32").Value = 4424.1665
$ws2.Range("J32").Value = 33226
$ws2.Range("K32").Value = 4424.1665
$ws2.Range("L32").Value = 33226
$ws2.Range("M32").Value = -4137.1665
$ws2.Range("N32").Value = -33800
# row 45 (diff @@ -9950)
$ws2.Range("H45").Value = 1800.6666
$ws2.Range("I45").Value = 3042.4
$ws2.Range("J45").Value = 1323.0769
$ws2.Range("K45").Value = 3042.4
$ws2.Range("L45").Value = 1323.0769
$ws2.Range("M45").Value = -2665.4
$ws2.Range("N45").Value = -2077.0769
# row 54 (diff @@ -10394)
$ws2.Range("H54").Value = 6049
$ws2.Range("I54").Value = 0
$ws2.Range("J54").Value = 6049
$ws2.Range("K54").Value = 0
$ws2.Range("L54").Value = 6049
$ws2.Range("N54").Value = -7587
# row 105 (diff @@ -12884)
$ws2.Range("H105").Value = 22370
$ws2.Range("I105").Value = 0
$ws2.Range("J105").Value = 22370
$ws2.Range("K105").Value = 0
$ws2.Range("L105").Value = 22370
$ws2.Range("N105").Value = -29358
# row 110 (diff @@ -13123)
$ws2.Range("H110").Value = 2357.3635
$ws2.Range("I110").Value = 2152.889
$ws2.Range("J110").Value = 3277.5
$ws2.Range("K110").Value = 2152.889
$ws2.Range("L110").Value = 3277.5
$ws2.Range("M110").Value = -107.8890000000001
$ws2.Range("N110").Value = -7367.5
# row 116 (diff @@ -13417)
$ws2.Range("H116").Value = 630.30304
$ws2.Range("I116").Value = 550.3333
$ws2.Range("J116").Value = 696.94446
$ws2.Range("K116").Value = 550.3333
$ws2.Range("L116").Value = 696.94446
$ws2.Range("M116").Value = 1743.6667
$ws2.Range("N116").Value = -5284.94446
# row 125 (diff @@ -13864)
$ws2.Range("H125").Value = 41703
$ws2.Range("I125").Value = 0
$ws2.Range("J125").Value = 41703
$ws2.Range("K125").Value = 0
$ws2.Range("L125").Value = 41703
$ws2.Range("N125").Value = -51543

$ws3 = $wb.Worksheets.Item("BSM")
# row 3 (diff @@ -14843)
$ws3.Range("H3").Value = 630.30304
$ws3.Range("I3").Value = 550.3333
$ws3.Range("J3").Value = 696.94446
$ws3.Range("K3").Value = 550.3333
$ws3.Range("L3").Value = 696.94446
$ws3.Range("M3").Value = -436.3333
$ws3.Range("N3").Value = -924.94446

$ws4 = $wb.Worksheets.Item("CRP")
# row 31 (diff @@ -23217)
$ws4.Range("H31").Value = 8335196
$ws4.Range("I31").Value = 13794129
$ws4.Range("J31").Value = 3140.2104
$ws4.Range("K31").Value = 13794129
$ws4.Range("L31").Value = 3140.2104
$ws4.Range("M31").Value = -13793834
$ws4.Range("N31").Value = -3730.2104
# row 34 (diff @@ -23370)
$ws4.Range("H34").Value = 8335196
$ws4.Range("I34").Value = 13794129
$ws4.Range("J34").Value = 3140.2104
$ws4.Range("K34").Value = 13794129
$ws4.Range("L34").Value = 3140.2104
$ws4.Range("M34").Value = -13793927
$ws4.Range("N34").Value = -3544.2104
# row 121 (diff @@ -27588)
$ws4.Range("H121").Value = 20615.072
$ws4.Range("I121").Value = 0
$ws4.Range("J121").Value = 20615.072
$ws4.Range("K121").Value = 0
$ws4.Range("L121").Value = 20615.072
$ws4.Range("N121").Value = -23235.072
# row 123 (diff @@ -27686)
$ws4.Range("H123").Value = 28000
$ws4.Range("I123").Value = 0
$ws4.Range("J123").Value = 28000
$ws4.Range("K123").Value = 0
$ws4.Range("L123").Value = 28000
$ws4.Range("N123").Value = -37800

$ws5 = $wb.Worksheets.Item("CUL")
# row 69 (diff @@ -32093)
$ws5.Range("H69").Value = 7966.4443
$ws5.Range("I69").Value = 4500
$ws5.Range("J69").Value = 8956.857
$ws5.Range("K69").Value = 13500
$ws5.Range("L69").Value = 26870.571
$ws5.Range("M69").Value = -12689
$ws5.Range("N69").Value = -28492.571
# row 72 (diff @@ -32249)
$ws5.Range("H72").Value = 7966.4443
$ws5.Range("I72").Value = 4500
$ws5.Range("J72").Value = 8956.857
$ws5.Range("K72").Value = 40500
$ws5.Range("L72").Value = 80611.713
$ws5.Range("M72").Value = -36444
$ws5.Range("N72").Value = -88723.713

$ws6 = $wb.Worksheets.Item("GSM")
# row 102 (diff @@ -40805)
$ws6.Range("H102").Value = 1100
$ws6.Range("I102").Value = 0
$ws6.Range("J102").Value = 1100
$ws6.Range("K102").Value = 0
$ws6.Range("L102").Value = 1100
$ws6.Range("M102").ClearContents()
$ws6.Range("N102").Value = -4344
# row 123 (diff @@ -41828)
$ws6.Range("H123").Value = 22133
$ws6.Range("I123").Value = 0
$ws6.Range("J123").Value = 22133
$ws6.Range("K123").Value = 0
$ws6.Range("L123").Value = 22133
$ws6.Range("N123").Value = -27033
# row 136 (diff @@ -42465)
$ws6.Range("H136").Value = 19907.482
$ws6.Range("I136").Value = 0
$ws6.Range("J136").Value = 19907.482
$ws6.Range("K136").Value = 0
$ws6.Range("L136").Value = 59722.446
$ws6.Range("N136").Value = -64822.446

$ws7 = $wb.Worksheets.Item("LTW")
# row 7 (diff @@ -43089)
$ws7.Range("H7").Value = 1740.2
$ws7.Range("I7").Value = 1008.7778
$ws7.Range("J7").Value = 2338.6365
$ws7.Range("K7").Value = 1008.7778
$ws7.Range("L7").Value = 2338.6365
$ws7.Range("M7").Value = -896.7778
$ws7.Range("N7").Value = -2562.6365
# row 26 (diff @@ -44014)
$ws7.Range("H26").Value = 11249.25
$ws7.Range("I26").Value = 5500
$ws7.Range("J26").Value = 16998.5
$ws7.Range("K26").Value = 5500
$ws7.Range("L26").Value = 16998.5
$ws7.Range("M26").Value = -5205
$ws7.Range("N26").Value = -17588.5
# row 31 (diff @@ -44259)
$ws7.Range("H31").Value = 952.8333
$ws7.Range("I31").Value = 400
$ws7.Range("J31").Value = 1229.25
$ws7.Range("K31").Value = 400
$ws7.Range("L31").Value = 1229.25
$ws7.Range("M31").Value = -152
$ws7.Range("N31").Value = -1725.25
# row 40 (diff @@ -44712)
$ws7.Range("H40").Value = 2286.8
$ws7.Range("I40").Value = 2374.2222
$ws7.Range("J40").Value = 1500
$ws7.Range("K40").Value = 2374.2222
$ws7.Range("L40").Value = 1500
$ws7.Range("M40").Value = -2238.2222
$ws7.Range("N40").Value = -1772
# row 46 (diff @@ -45006)
$ws7.Range("H46").Value = 990.4
$ws7.Range("I46").Value = 851.3333
$ws7.Range("J46").Value = 1050
$ws7.Range("K46").Value = 851.3333
$ws7.Range("L46").Value = 1050
$ws7.Range("M46").Value = -663.3333
$ws7.Range("N46").Value = -1426
# row 122 (diff @@ -48703)
$ws7.Range("H122").Value = 5117.4375
$ws7.Range("I122").Value = 5732.6665
$ws7.Range("J122").Value = 4326.4287
$ws7.Range("K122").Value = 17197.9995
$ws7.Range("L122").Value = 12979.2861
$ws7.Range("M122").Value = -14747.9995
$ws7.Range("N122").Value = -17879.2861
# row 126 (diff @@ -48905)
$ws7.Range("H126").Value = 1740.2
$ws7.Range("I126").Value = 1008.7778
$ws7.Range("J126").Value = 2338.6365
$ws7.Range("K126").Value = 3026.3334
$ws7.Range("L126").Value = 7015.9095
$ws7.Range("M126").Value = -556.3334
$ws7.Range("N126").Value = -11955.9095

$ws8 = $wb.Worksheets.Item("WVR")
# row 37 (diff @@ -51465)
$ws8.Range("H37").Value = 3000
$ws8.Range("I37").Value = 3000
$ws8.Range("J37").Value = 0
$ws8.Range("K37").Value = 3000
$ws8.Range("L37").Value = 0
$ws8.Range("M37").Value = -2797
# row 122 (diff @@ -55612)
$ws8.Range("H122").Value = 1405.6428
$ws8.Range("I122").Value = 969.2857
$ws8.Range("J122").Value = 1842
$ws8.Range("K122").Value = 2907.8571
$ws8.Range("L122").Value = 5526
$ws8.Range("M122").Value = -457.8571000000002
$ws8.Range("N122").Value = -10426
# row 126 (diff @@ -55811)
$ws8.Range("H126").Value = 3395.8667
$ws8.Range("I126").Value = 3911.5
$ws8.Range("J126").Value = 1333.3334
$ws8.Range("K126").Value = 11734.5
$ws8.Range("L126").Value = 4000.0002
$ws8.Range("M126").Value = -9264.5
$ws8.Range("N126").Value = -8940.0002
